$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 37: fill in the new "16/6/2025" weekly entry ---
$ws.Range("D37").Value = "16/6/2025"
$ws.Range("E37").Value = 408
$ws.Range("F37").Value = 518
$ws.Range("G37").Value = 0
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 1012
$ws.Range("J37").Value = "N/A"

# --- Column D formatting: right-align the date cells in rows 33-44 ---
# (matches the rest of the table's date column style)
$ws.Range("D33:D44").HorizontalAlignment = -4152

# --- View state: scroll down, zoom in, move the active selection ---
$win = $excel.ActiveWindow
$win.Zoom = 85
$win.ScrollRow = 17
$win.ScrollColumn = 1
$ws.Range("D40").Select()
